$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -9.024199999999999
$ws.Range("D9").Value = -8.683999999999992
$ws.Range("B11").Value = 6.834700000000002
$ws.Range("C11").Value = -11.7758
$ws.Range("B12").Value = 4.674599999999997
$ws.Range("D13").Value = -8.313299999999996
$ws.Range("D14").Value = -8.385799999999996
$ws.Range("B15").Value = 5.826800000000002
$ws.Range("D19").Value = -8.625099999999998
$ws.Range("D21").Value = -8.258200000000002
$ws.Range("D22").Value = -7.508499999999997
$ws.Range("C23").Value = -12.08940000000001
$ws.Range("D24").Value = -7.713199999999995
$ws.Range("D26").Value = -7.870999999999998
$ws.Range("B27").Value = 5.707500000000001
$ws.Range("B28").Value = 6.105599999999996
$ws.Range("C28").Value = -13.0577
$ws.Range("B31").Value = 4.661199999999997
$ws.Range("B32").Value = 6.659899999999999
$ws.Range("C32").Value = -11.4788
$ws.Range("C34").Value = -11.88280000000001
$ws.Range("B36").Value = 9.522100000000002
$ws.Range("C36").Value = -11.93770000000001
$ws.Range("C37").Value = -13.33600000000001
$ws.Range("B38").Value = 6.446100000000001
$ws.Range("D38").Value = -7.624100000000003
$ws.Range("D41").Value = -8.142899999999999
$ws.Range("C42").Value = -12.49110000000001
$ws.Range("B46").Value = 6.217099999999999
$ws.Range("C49").Value = -13.40049999999999
$ws.Range("D52").Value = -7.811900000000003
$ws.Range("B54").Value = 5.197400000000004
$ws.Range("C54").Value = -14.0005
$ws.Range("B55").Value = 5.184999999999998
$ws.Range("B56").Value = 5.829800000000001
$ws.Range("D56").Value = -7.902600000000002
$ws.Range("B67").Value = 5.403599999999996
$ws.Range("B69").Value = 5.343899999999996
$ws.Range("D71").Value = -7.208599999999994
$ws.Range("B72").Value = 6.498299999999996
$ws.Range("D72").Value = -7.080099999999998
$ws.Range("B73").Value = 9.256100000000004
$ws.Range("C78").Value = -12.95820000000001
$ws.Range("D78").Value = -8.310500000000001
$ws.Range("C80").Value = -11.81530000000001
$ws.Range("B83").Value = 6.013299999999997
$ws.Range("D83").Value = -8.8398
$ws.Range("D85").Value = -8.966899999999999
$ws.Range("B86").Value = 5.869000000000002
$ws.Range("D86").Value = -7.957900000000003
$ws.Range("D90").Value = -7.108199999999991
$ws.Range("B91").Value = 4.842799999999999
$ws.Range("B93").Value = 5.431200000000005
$ws.Range("D96").Value = -8.385299999999992
$ws.Range("C97").Value = -11.5237
$ws.Range("B99").Value = 6.413800000000001
$ws.Range("C99").Value = -11.99570000000001
$ws.Range("C100").Value = -11.9709
$ws.Range("C101").Value = -13.1865
$ws.Range("D103").Value = -8.4125
$ws.Range("B104").Value = 9.3932
$ws.Range("B105").Value = 8.488700000000001
